$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 54, pushing existing rows 54:78 down to 55:79
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the weekly Jengibre price data
$ws.Cells.Item(54, 1).Value = 8
$ws.Cells.Item(54, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44879
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(54, 6).Value = 100114007
$ws.Cells.Item(54, 7).Value = "Jengibre"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 400
$ws.Cells.Item(54, 11).Value = 13500
$ws.Cells.Item(54, 12).Value = 14000
$ws.Cells.Item(54, 13).Value = 13750
$ws.Cells.Item(54, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(54, 15).Value = "Perú"
$ws.Cells.Item(54, 16).Value = 1058
$ws.Cells.Item(54, 17).Value = 13
$ws.Cells.Item(54, 18).Value = "Hortaliza"
